# Corrected excel sheets for application fix issues
#
# 1) "Summary" sheet: widen the held selection from A7:XFD12 to A7:XFD14
# 2) "Repayment schedule" sheet: move the held selection from B7 to A9:XFD9
# 3) "Transactions" sheet: correct transaction IDs in A2 (18 -> 67) and A3 (17 -> 66)

$wb = $excel.ActiveWorkbook

# --- Summary sheet: update selection ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
$wsSummary.Range("A7:XFD14").Select()

# --- Repayment schedule sheet: update selection ---
$wsRepayment = $wb.Worksheets.Item("Repayment schedule")
$wsRepayment.Activate()
$wsRepayment.Range("A9:XFD9").Select()

# --- Transactions sheet: correct the ID values, keep it the active sheet ---
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Activate()
$wsTransactions.Range("A2").Value = 67
$wsTransactions.Range("A3").Value = 66
